# Natmi following Dr Hou advice
# Recompute the LR-pairs (Vegfa-Flt1) table: existing rows 2-13 get refreshed
# metric values (ligand/receptor expressing-cell counts, expression values,
# specificities, edge weights), and four new rows (14-17) are appended for
# the "sCs" sending cluster against each of the four target clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vegfa"
$ws.Cells.Item(2, 3).Value = "Flt1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 28.479168
$ws.Cells.Item(2, 8).Value = 85.437504
$ws.Cells.Item(2, 9).Value = 0.4446244458164738
$ws.Cells.Item(2, 10).Value = 0.4446244458164738
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 100.95625
$ws.Cells.Item(2, 14).Value = 302.86875
$ws.Cells.Item(2, 15).Value = 0.9361599519103526
$ws.Cells.Item(2, 16).Value = 0.9361599519103525
$ws.Cells.Item(2, 17).Value = 2875.1500044
$ws.Cells.Item(2, 18).Value = 25876.3500396
$ws.Cells.Item(2, 19).Value = 0.4162395998137173
$ws.Cells.Item(2, 20).Value = 0.4162395998137173

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vegfa"
$ws.Cells.Item(3, 3).Value = "Flt1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 28.479168
$ws.Cells.Item(3, 8).Value = 85.437504
$ws.Cells.Item(3, 9).Value = 0.4446244458164738
$ws.Cells.Item(3, 10).Value = 0.4446244458164738
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.3890956666666667
$ws.Cells.Item(3, 14).Value = 1.167287
$ws.Cells.Item(3, 15).Value = 0.003608055772626195
$ws.Cells.Item(3, 16).Value = 0.003608055772626195
$ws.Cells.Item(3, 17).Value = 11.081120859072
$ws.Cells.Item(3, 18).Value = 99.730087731648
$ws.Cells.Item(3, 19).Value = 0.001604229798378851
$ws.Cells.Item(3, 20).Value = 0.001604229798378851

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vegfa"
$ws.Cells.Item(4, 3).Value = "Flt1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 28.479168
$ws.Cells.Item(4, 8).Value = 85.437504
$ws.Cells.Item(4, 9).Value = 0.4446244458164738
$ws.Cells.Item(4, 10).Value = 0.4446244458164738
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.996608333333334
$ws.Cells.Item(4, 14).Value = 8.989825000000002
$ws.Cells.Item(4, 15).Value = 0.02778733078167519
$ws.Cells.Item(4, 16).Value = 0.02778733078167519
$ws.Cells.Item(4, 17).Value = 85.34091215520002
$ws.Cells.Item(4, 18).Value = 768.0682093968002
$ws.Cells.Item(4, 19).Value = 0.01235492654952138
$ws.Cells.Item(4, 20).Value = 0.01235492654952138

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Vegfa"
$ws.Cells.Item(5, 3).Value = "Flt1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 28.479168
$ws.Cells.Item(5, 8).Value = 85.437504
$ws.Cells.Item(5, 9).Value = 0.4446244458164738
$ws.Cells.Item(5, 10).Value = 0.4446244458164738
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.498858666666667
$ws.Cells.Item(5, 14).Value = 10.496576
$ws.Cells.Item(5, 15).Value = 0.03244466153534613
$ws.Cells.Item(5, 16).Value = 0.03244466153534613
$ws.Cells.Item(5, 17).Value = 99.64458377625601
$ws.Cells.Item(5, 18).Value = 896.8012539863041
$ws.Cells.Item(5, 19).Value = 0.01442568965485634
$ws.Cells.Item(5, 20).Value = 0.01442568965485634

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vegfa"
$ws.Cells.Item(6, 3).Value = "Flt1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 18.12667766666667
$ws.Cells.Item(6, 8).Value = 54.380033
$ws.Cells.Item(6, 9).Value = 0.2829985767855128
$ws.Cells.Item(6, 10).Value = 0.2829985767855128
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 100.95625
$ws.Cells.Item(6, 14).Value = 302.86875
$ws.Cells.Item(6, 15).Value = 0.9361599519103526
$ws.Cells.Item(6, 16).Value = 0.9361599519103525
$ws.Cells.Item(6, 17).Value = 1830.001402185417
$ws.Cells.Item(6, 18).Value = 16470.01261966875
$ws.Cells.Item(6, 19).Value = 0.2649319340342239
$ws.Cells.Item(6, 20).Value = 0.2649319340342239

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vegfa"
$ws.Cells.Item(7, 3).Value = "Flt1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 18.12667766666667
$ws.Cells.Item(7, 8).Value = 54.380033
$ws.Cells.Item(7, 9).Value = 0.2829985767855128
$ws.Cells.Item(7, 10).Value = 0.2829985767855128
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.3890956666666667
$ws.Cells.Item(7, 14).Value = 1.167287
$ws.Cells.Item(7, 15).Value = 0.003608055772626195
$ws.Cells.Item(7, 16).Value = 0.003608055772626195
$ws.Cells.Item(7, 17).Value = 7.053011731163446
$ws.Cells.Item(7, 18).Value = 63.477105580471
$ws.Cells.Item(7, 19).Value = 0.001021074648615967
$ws.Cells.Item(7, 20).Value = 0.001021074648615967

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Vegfa"
$ws.Cells.Item(8, 3).Value = "Flt1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 18.12667766666667
$ws.Cells.Item(8, 8).Value = 54.380033
$ws.Cells.Item(8, 9).Value = 0.2829985767855128
$ws.Cells.Item(8, 10).Value = 0.2829985767855128
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.996608333333334
$ws.Cells.Item(8, 14).Value = 8.989825000000002
$ws.Cells.Item(8, 15).Value = 0.02778733078167519
$ws.Cells.Item(8, 16).Value = 0.02778733078167519
$ws.Cells.Item(8, 17).Value = 54.31855335158058
$ws.Cells.Item(8, 18).Value = 488.8669801642251
$ws.Cells.Item(8, 19).Value = 0.00786377506388235
$ws.Cells.Item(8, 20).Value = 0.00786377506388235

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Vegfa"
$ws.Cells.Item(9, 3).Value = "Flt1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 18.12667766666667
$ws.Cells.Item(9, 8).Value = 54.380033
$ws.Cells.Item(9, 9).Value = 0.2829985767855128
$ws.Cells.Item(9, 10).Value = 0.2829985767855128
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.498858666666667
$ws.Cells.Item(9, 14).Value = 10.496576
$ws.Cells.Item(9, 15).Value = 0.03244466153534613
$ws.Cells.Item(9, 16).Value = 0.03244466153534613
$ws.Cells.Item(9, 17).Value = 63.42268325188979
$ws.Cells.Item(9, 18).Value = 570.8041492670081
$ws.Cells.Item(9, 19).Value = 0.009181793038790626
$ws.Cells.Item(9, 20).Value = 0.009181793038790626

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Vegfa"
$ws.Cells.Item(10, 3).Value = "Flt1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 11.513346
$ws.Cells.Item(10, 8).Value = 34.540038
$ws.Cells.Item(10, 9).Value = 0.179749460544048
$ws.Cells.Item(10, 10).Value = 0.179749460544048
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 100.95625
$ws.Cells.Item(10, 14).Value = 302.86875
$ws.Cells.Item(10, 15).Value = 0.9361599519103526
$ws.Cells.Item(10, 16).Value = 0.9361599519103525
$ws.Cells.Item(10, 17).Value = 1162.3442371125
$ws.Cells.Item(10, 18).Value = 10461.0981340125
$ws.Cells.Item(10, 19).Value = 0.1682742463388278
$ws.Cells.Item(10, 20).Value = 0.1682742463388278

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Vegfa"
$ws.Cells.Item(11, 3).Value = "Flt1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 11.513346
$ws.Cells.Item(11, 8).Value = 34.540038
$ws.Cells.Item(11, 9).Value = 0.179749460544048
$ws.Cells.Item(11, 10).Value = 0.179749460544048
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3890956666666667
$ws.Cells.Item(11, 14).Value = 1.167287
$ws.Cells.Item(11, 15).Value = 0.003608055772626195
$ws.Cells.Item(11, 16).Value = 0.003608055772626195
$ws.Cells.Item(11, 17).Value = 4.479793037434
$ws.Cells.Item(11, 18).Value = 40.318137336906
$ws.Cells.Item(11, 19).Value = 0.0006485460787423969
$ws.Cells.Item(11, 20).Value = 0.0006485460787423969

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Vegfa"
$ws.Cells.Item(12, 3).Value = "Flt1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 11.513346
$ws.Cells.Item(12, 8).Value = 34.540038
$ws.Cells.Item(12, 9).Value = 0.179749460544048
$ws.Cells.Item(12, 10).Value = 0.179749460544048
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.996608333333334
$ws.Cells.Item(12, 14).Value = 8.989825000000002
$ws.Cells.Item(12, 15).Value = 0.02778733078167519
$ws.Cells.Item(12, 16).Value = 0.02778733078167519
$ws.Cells.Item(12, 17).Value = 34.50098856815001
$ws.Cells.Item(12, 18).Value = 310.5088971133501
$ws.Cells.Item(12, 19).Value = 0.004994757717965136
$ws.Cells.Item(12, 20).Value = 0.004994757717965136

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Vegfa"
$ws.Cells.Item(13, 3).Value = "Flt1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 11.513346
$ws.Cells.Item(13, 8).Value = 34.540038
$ws.Cells.Item(13, 9).Value = 0.179749460544048
$ws.Cells.Item(13, 10).Value = 0.179749460544048
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.498858666666667
$ws.Cells.Item(13, 14).Value = 10.496576
$ws.Cells.Item(13, 15).Value = 0.03244466153534613
$ws.Cells.Item(13, 16).Value = 0.03244466153534613
$ws.Cells.Item(13, 17).Value = 40.28357043443201
$ws.Cells.Item(13, 18).Value = 362.5521339098881
$ws.Cells.Item(13, 19).Value = 0.005831910408512692
$ws.Cells.Item(13, 20).Value = 0.005831910408512692

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Vegfa"
$ws.Cells.Item(14, 3).Value = "Flt1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 5.932994999999999
$ws.Cells.Item(14, 8).Value = 17.798985
$ws.Cells.Item(14, 9).Value = 0.09262751685396531
$ws.Cells.Item(14, 10).Value = 0.09262751685396531
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 100.95625
$ws.Cells.Item(14, 14).Value = 302.86875
$ws.Cells.Item(14, 15).Value = 0.9361599519103526
$ws.Cells.Item(14, 16).Value = 0.9361599519103525
$ws.Cells.Item(14, 17).Value = 598.9729264687498
$ws.Cells.Item(14, 18).Value = 5390.756338218749
$ws.Cells.Item(14, 19).Value = 0.08671417172358353
$ws.Cells.Item(14, 20).Value = 0.08671417172358353

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Vegfa"
$ws.Cells.Item(15, 3).Value = "Flt1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 5.932994999999999
$ws.Cells.Item(15, 8).Value = 17.798985
$ws.Cells.Item(15, 9).Value = 0.09262751685396531
$ws.Cells.Item(15, 10).Value = 0.09262751685396531
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.3890956666666667
$ws.Cells.Item(15, 14).Value = 1.167287
$ws.Cells.Item(15, 15).Value = 0.003608055772626195
$ws.Cells.Item(15, 16).Value = 0.003608055772626195
$ws.Cells.Item(15, 17).Value = 2.308502644855
$ws.Cells.Item(15, 18).Value = 20.776523803695
$ws.Cells.Item(15, 19).Value = 0.0003342052468889797
$ws.Cells.Item(15, 20).Value = 0.0003342052468889797

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Vegfa"
$ws.Cells.Item(16, 3).Value = "Flt1"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 5.932994999999999
$ws.Cells.Item(16, 8).Value = 17.798985
$ws.Cells.Item(16, 9).Value = 0.09262751685396531
$ws.Cells.Item(16, 10).Value = 0.09262751685396531
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 2.996608333333334
$ws.Cells.Item(16, 14).Value = 8.989825000000002
$ws.Cells.Item(16, 15).Value = 0.02778733078167519
$ws.Cells.Item(16, 16).Value = 0.02778733078167519
$ws.Cells.Item(16, 17).Value = 17.778862258625
$ws.Cells.Item(16, 18).Value = 160.009760327625
$ws.Cells.Item(16, 19).Value = 0.002573871450306328
$ws.Cells.Item(16, 20).Value = 0.002573871450306328

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Vegfa"
$ws.Cells.Item(17, 3).Value = "Flt1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 5.932994999999999
$ws.Cells.Item(17, 8).Value = 17.798985
$ws.Cells.Item(17, 9).Value = 0.09262751685396531
$ws.Cells.Item(17, 10).Value = 0.09262751685396531
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 3.498858666666667
$ws.Cells.Item(17, 14).Value = 10.496576
$ws.Cells.Item(17, 15).Value = 0.03244466153534613
$ws.Cells.Item(17, 16).Value = 0.03244466153534613
$ws.Cells.Item(17, 17).Value = 20.75871097504
$ws.Cells.Item(17, 18).Value = 186.82839877536
$ws.Cells.Item(17, 19).Value = 0.003005268433186473
$ws.Cells.Item(17, 20).Value = 0.003005268433186473
